# Update Betfair back/lay odds values on the active (only) worksheet.
# These are straightforward numeric re-quotes for existing rows; no rows/
# columns are added or removed, and no formatting changes are required.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = 2.88
$ws.Range("H3").Value = 2.2
$ws.Range("K3").Value = 5.2
$ws.Range("N3").Value = 3.2
$ws.Range("Q3").Value = 1.91
$ws.Range("S3").Value = 3.25

$ws.Range("G4").Value = 1.92
$ws.Range("H4").Value = 2.08

$ws.Range("H5").Value = 2.46
$ws.Range("J5").Value = 2.7
$ws.Range("K5").Value = 7.2
$ws.Range("Q5").Value = 1.41

$ws.Range("F7").Value = 5.6
$ws.Range("G7").Value = 240
$ws.Range("H7").Value = 1.32
$ws.Range("K7").Value = 6.8
$ws.Range("P7").Value = 1.92
$ws.Range("Q7").Value = 1.86

$ws.Range("J8").Value = 3.25
$ws.Range("P8").Value = 1.76

$ws.Range("F9").Value = 2.82
$ws.Range("G9").Value = 3.9
$ws.Range("H9").Value = 2.14
$ws.Range("I9").Value = 2.76
$ws.Range("J9").Value = 2.62
$ws.Range("K9").Value = 5.8
$ws.Range("P9").Value = 2.06

$ws.Range("N10").Value = 1.71
$ws.Range("P10").Value = 1.72
$ws.Range("S10").Value = 3.3

$ws.Range("P11").Value = 3.1

$ws.Range("G12").Value = 3.65
$ws.Range("H12").Value = 2.14
$ws.Range("J12").Value = 2.72
$ws.Range("K12").Value = 6.8

$ws.Range("H13").Value = 1.87
$ws.Range("K13").Value = 8.199999999999999
$ws.Range("P13").Value = 3.1

$ws.Range("Y16").Value = 23
$ws.Range("Z16").Value = 85
$ws.Range("AA16").Value = 470
$ws.Range("AD16").Value = 36

$ws.Range("F17").Value = 2.42
$ws.Range("G17").Value = 2.64
$ws.Range("H17").Value = 3.55
$ws.Range("I17").Value = 4.1
$ws.Range("J17").Value = 2.84
$ws.Range("K17").Value = 3.2
$ws.Range("Q17").Value = 2.8
$ws.Range("W17").Value = 1.61

$ws.Range("N18").Value = 3.3

$ws.Range("F19").Value = 6.6
$ws.Range("J19").Value = 4.3
$ws.Range("P19").Value = 2.04

$ws.Range("G20").Value = 2.22
$ws.Range("H20").Value = 3.6
$ws.Range("I20").Value = 4
$ws.Range("U20").Value = 2.16
$ws.Range("W20").Value = 1.81
$ws.Range("AN20").Value = 18.5
